$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Move the UAV id values up one row in column C:
# C4 gets the value previously in C5 ("sig kadet")
# C5 gets the value previously in C6 ("valiant")
# C6 is cleared
$ws.Range("C4").Value = "sig kadet"
$ws.Range("C5").Value = "valiant"
$ws.Range("C6").ClearContents()

# Update the active cell selection to D11 as recorded in the saved file
$ws.Range("D11").Select()
